# TDE01 - apply Arial/12pt formatting to every paragraph and split the
# trailing "." of the short answer paragraphs (01, 02, 04, 11) into its
# own run, matching the target revision.

$d = $word.ActiveDocument

function Format-Range($rng) {
    $rng.Font.Name = "Arial"
    $rng.Font.NameBi = "Arial"
    $rng.Font.Size = 12
    $rng.Font.SizeBi = 12
}

# Paragraphs whose text needs a literal "." appended and split into its
# own run (they currently end with just the letter, e.g. "01 - A").
$paragraphsNeedingPeriod = @(1, 2, 4)

# Paragraph whose text already ends with "." but still needs that "."
# pulled out into a separate run.
$paragraphsNeedingSplitOnly = @(11)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {

    $p = $d.Paragraphs($i)

    # Format the whole paragraph (sets pPr/rPr mark formatting plus the
    # existing run(s)).
    Format-Range $p.Range

    if ($paragraphsNeedingPeriod -contains $i) {
        # Insert the missing "." right before the paragraph mark.
        $insertPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)
        $insertPoint.InsertAfter(".")

        # Re-apply formatting to just the new "." so it lands in its own run.
        $periodRange = $d.Range($p.Range.End - 2, $p.Range.End - 1)
        Format-Range $periodRange
    }
    elseif ($paragraphsNeedingSplitOnly -contains $i) {
        # The "." is already part of the text; force it into its own run
        # by nudging its formatting (same value in, same value out).
        $periodRange = $d.Range($p.Range.End - 2, $p.Range.End - 1)
        $periodRange.Font.Size = 13
        $periodRange.Font.Size = 12
    }
}
